$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("spellDictionary")

# Add three new spell animation names (shared strings will be created automatically).
# Set them in the same order they first appear as new unique strings so the
# shared string table ends up ordered: anim_spell_ice, anim_spell_empower, anim_spell_heal

# Row 26: cryo / ice -> anim_spell_ice (new unique string #1)
$ws.Range("D26").Value = "anim_spell_ice"

# Row 30: aimed -> anim_spell_empower (new unique string #2)
$ws.Range("D30").Value = "anim_spell_empower"

# Row 31: empowered -> anim_spell_empower
$ws.Range("D31").Value = "anim_spell_empower"

# Row 32: widened -> anim_spell_empower
$ws.Range("D32").Value = "anim_spell_empower"

# Row 33: shank -> anim_spell_empower
$ws.Range("D33").Value = "anim_spell_empower"

# Row 34: miss -> anim_spell_empower
$ws.Range("D34").Value = "anim_spell_empower"

# Row 5: selfcare / heal self -> anim_spell_heal (new unique string #3)
$ws.Range("D5").Value = "anim_spell_heal"

# Column D width change (~14.66 characters wide)
$ws.Columns.Item(4).ColumnWidth = 13.83

# Sheet view changes: remove topLeftCell scroll position, change selection to D6
$ws.Activate()
$ws.Range("D6").Select()
